$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.950239
$ws.Range("H2").Value = 152.850717
$ws.Range("I2").Value = 0.1520006117784607
$ws.Range("J2").Value = 0.1540898474582185
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.99123566666667
$ws.Range("N2").Value = 32.973707
$ws.Range("O2").Value = 0.06247112414817082
$ws.Range("P2").Value = 0.06381790718517791
$ws.Range("Q2").Value = 560.0060841219911
$ws.Range("R2").Value = 5040.05475709792
$ws.Range("S2").Value = 0.009495649089010134
$ws.Range("T2").Value = 0.00983369158326681

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.950239
$ws.Range("H3").Value = 152.850717
$ws.Range("I3").Value = 0.1520006117784607
$ws.Range("J3").Value = 0.1540898474582185
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 16.58451233333333
$ws.Range("N3").Value = 49.753537
$ws.Range("O3").Value = 0.09426175184784683
$ws.Range("P3").Value = 0.09629389277948987
$ws.Range("Q3").Value = 844.984867081781
$ws.Range("R3").Value = 7604.86380373603
$ws.Range("S3").Value = 0.01432784394818217
$ws.Range("T3").Value = 0.01483791124954964

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.950239
$ws.Range("H4").Value = 152.850717
$ws.Range("I4").Value = 0.1520006117784607
$ws.Range("J4").Value = 0.1540898474582185
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 78.20455633333331
$ws.Range("N4").Value = 234.613669
$ws.Range("O4").Value = 0.4444929301687811
$ws.Range("P4").Value = 0.4540755260734272
$ws.Range("Q4").Value = 3984.540836072296
$ws.Range("R4").Value = 35860.86752465067
$ws.Range("S4").Value = 0.06756319731685535
$ws.Range("T4").Value = 0.0699684285471647

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.950239
$ws.Range("H5").Value = 152.850717
$ws.Range("I5").Value = 0.1520006117784607
$ws.Range("J5").Value = 0.1540898474582185
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 59.02182200000001
$ws.Range("N5").Value = 177.065466
$ws.Range("O5").Value = 0.3354636076811054
$ws.Range("P5").Value = 0.3426956961462742
$ws.Range("Q5").Value = 3007.175937115459
$ws.Range("R5").Value = 27064.58343403913
$ws.Range("S5").Value = 0.05099067359693756
$ws.Range("T5").Value = 0.05280592754376739

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 50.950239
$ws.Range("H6").Value = 152.850717
$ws.Range("I6").Value = 0.1520006117784607
$ws.Range("J6").Value = 0.1540898474582185
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 11.1389315
$ws.Range("N6").Value = 22.277863
$ws.Range("O6").Value = 0.06331058615409579
$ws.Range("P6").Value = 0.04311697781563077
$ws.Range("Q6").Value = 567.5312221296285
$ws.Range("R6").Value = 3405.187332777771
$ws.Range("S6").Value = 0.009623247827475504
$ws.Range("T6").Value = 0.006643888534469935

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.19715733333333
$ws.Range("H7").Value = 45.591472
$ws.Range("I7").Value = 0.04533790728558088
$ws.Range("J7").Value = 0.0459610730244441
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.99123566666667
$ws.Range("N7").Value = 32.973707
$ws.Range("O7").Value = 0.06247112414817082
$ws.Range("P7").Value = 0.06381790718517791
$ws.Range("Q7").Value = 167.0355377140782
$ws.Range("R7").Value = 1503.319839426704
$ws.Range("S7").Value = 0.002832310034655781
$ws.Range("T7").Value = 0.002933139492405157

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.19715733333333
$ws.Range("H8").Value = 45.591472
$ws.Range("I8").Value = 0.04533790728558088
$ws.Range("J8").Value = 0.0459610730244441
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.58451233333333
$ws.Range("N8").Value = 49.753537
$ws.Range("O8").Value = 0.09426175184784683
$ws.Range("P8").Value = 0.09629389277948987
$ws.Range("Q8").Value = 252.0374432262737
$ws.Range("R8").Value = 2268.336989036464
$ws.Range("S8").Value = 0.004273630565854111
$ws.Range("T8").Value = 0.004425770637846124

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.19715733333333
$ws.Range("H9").Value = 45.591472
$ws.Range("I9").Value = 0.04533790728558088
$ws.Range("J9").Value = 0.0459610730244441
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 78.20455633333331
$ws.Range("N9").Value = 234.613669
$ws.Range("O9").Value = 0.4444929301687811
$ws.Range("P9").Value = 0.4540755260734272
$ws.Range("Q9").Value = 1188.486946781196
$ws.Range("R9").Value = 10696.38252103077
$ws.Range("S9").Value = 0.02015237925708837
$ws.Range("T9").Value = 0.02086979841247365

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.19715733333333
$ws.Range("H10").Value = 45.591472
$ws.Range("I10").Value = 0.04533790728558088
$ws.Range("J10").Value = 0.0459610730244441
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 59.02182200000001
$ws.Range("N10").Value = 177.065466
$ws.Range("O10").Value = 0.3354636076811054
$ws.Range("P10").Value = 0.3426956961462742
$ws.Range("Q10").Value = 896.9639150339947
$ws.Range("R10").Value = 8072.675235305952
$ws.Range("S10").Value = 0.01520921794273243
$ws.Range("T10").Value = 0.01575066191574161

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 15.19715733333333
$ws.Range("H11").Value = 45.591472
$ws.Range("I11").Value = 0.04533790728558088
$ws.Range("J11").Value = 0.0459610730244441
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 11.1389315
$ws.Range("N11").Value = 22.277863
$ws.Range("O11").Value = 0.06331058615409579
$ws.Range("P11").Value = 0.04311697781563077
$ws.Range("Q11").Value = 169.2800945307226
$ws.Range("R11").Value = 1015.680567184336
$ws.Range("S11").Value = 0.002870369485250175
$ws.Range("T11").Value = 0.001981702565977542

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 124.2523
$ws.Range("H12").Value = 372.7569
$ws.Range("I12").Value = 0.3706837491945981
$ws.Range("J12").Value = 0.375778766284743
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.99123566666667
$ws.Range("N12").Value = 32.973707
$ws.Range("O12").Value = 0.06247112414817082
$ws.Range("P12").Value = 0.06381790718517791
$ws.Range("Q12").Value = 1365.686311425367
$ws.Range("R12").Value = 12291.1768028283
$ws.Range("S12").Value = 0.02315703051564515
$ws.Range("T12").Value = 0.02398141442892039

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 124.2523
$ws.Range("H13").Value = 372.7569
$ws.Range("I13").Value = 0.3706837491945981
$ws.Range("J13").Value = 0.375778766284743
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.58451233333333
$ws.Range("N13").Value = 49.753537
$ws.Range("O13").Value = 0.09426175184784683
$ws.Range("P13").Value = 0.09629389277948987
$ws.Range("Q13").Value = 2060.663801795033
$ws.Range("R13").Value = 18545.9742161553
$ws.Range("S13").Value = 0.0349412995806107
$ws.Range("T13").Value = 0.03618520022943203

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 124.2523
$ws.Range("H14").Value = 372.7569
$ws.Range("I14").Value = 0.3706837491945981
$ws.Range("J14").Value = 0.375778766284743
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 78.20455633333331
$ws.Range("N14").Value = 234.613669
$ws.Range("O14").Value = 0.4444929301687811
$ws.Range("P14").Value = 0.4540755260734272
$ws.Range("Q14").Value = 9717.09599489623
$ws.Range("R14").Value = 87453.86395406608
$ws.Range("S14").Value = 0.1647663058454565
$ws.Range("T14").Value = 0.1706319409879681

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 124.2523
$ws.Range("H15").Value = 372.7569
$ws.Range("I15").Value = 0.3706837491945981
$ws.Range("J15").Value = 0.375778766284743
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.02182200000001
$ws.Range("N15").Value = 177.065466
$ws.Range("O15").Value = 0.3354636076811054
$ws.Range("P15").Value = 0.3426956961462742
$ws.Range("Q15").Value = 7333.5971336906
$ws.Range("R15").Value = 66002.3742032154
$ws.Range("S15").Value = 0.1243509078135779
$ws.Range("T15").Value = 0.1287777659089381

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 124.2523
$ws.Range("H16").Value = 372.7569
$ws.Range("I16").Value = 0.3706837491945981
$ws.Range("J16").Value = 0.375778766284743
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 11.1389315
$ws.Range("N16").Value = 22.277863
$ws.Range("O16").Value = 0.06331058615409579
$ws.Range("P16").Value = 0.04311697781563077
$ws.Range("Q16").Value = 1384.03785841745
$ws.Range("R16").Value = 8304.227150504699
$ws.Range("S16").Value = 0.02346820543930784
$ws.Range("T16").Value = 0.01620244472948437

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 131.1635103333333
$ws.Range("H17").Value = 393.490531
$ws.Range("I17").Value = 0.3913020665845575
$ws.Range("J17").Value = 0.3966804807205673
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 10.99123566666667
$ws.Range("N17").Value = 32.973707
$ws.Range("O17").Value = 0.06247112414817082
$ws.Range("P17").Value = 0.06381790718517791
$ws.Range("Q17").Value = 1441.649052940935
$ws.Range("R17").Value = 12974.84147646842
$ws.Range("S17").Value = 0.0244450799810397
$ws.Range("T17").Value = 0.02531531810079692

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 131.1635103333333
$ws.Range("H18").Value = 393.490531
$ws.Range("I18").Value = 0.3913020665845575
$ws.Range("J18").Value = 0.3966804807205673
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 16.58451233333333
$ws.Range("N18").Value = 49.753537
$ws.Range("O18").Value = 0.09426175184784683
$ws.Range("P18").Value = 0.09629389277948987
$ws.Range("Q18").Value = 2175.282854806461
$ws.Range("R18").Value = 19577.54569325815
$ws.Range("S18").Value = 0.0368848182979432
$ws.Range("T18").Value = 0.03819790767822281

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 131.1635103333333
$ws.Range("H19").Value = 393.490531
$ws.Range("I19").Value = 0.3913020665845575
$ws.Range("J19").Value = 0.3966804807205673
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 78.20455633333331
$ws.Range("N19").Value = 234.613669
$ws.Range("O19").Value = 0.4444929301687811
$ws.Range("P19").Value = 0.4540755260734272
$ws.Range("Q19").Value = 10257.58413274091
$ws.Range("R19").Value = 92318.25719466823
$ws.Range("S19").Value = 0.1739310021572694
$ws.Range("T19").Value = 0.1801228979662516

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 131.1635103333333
$ws.Range("H20").Value = 393.490531
$ws.Range("I20").Value = 0.3913020665845575
$ws.Range("J20").Value = 0.3966804807205673
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 59.02182200000001
$ws.Range("N20").Value = 177.065466
$ws.Range("O20").Value = 0.3354636076811054
$ws.Range("P20").Value = 0.3426956961462742
$ws.Range("Q20").Value = 7741.509359789162
$ws.Range("R20").Value = 69673.58423810246
$ws.Range("S20").Value = 0.1312676029495278
$ws.Range("T20").Value = 0.1359406934881735

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 131.1635103333333
$ws.Range("H21").Value = 393.490531
$ws.Range("I21").Value = 0.3913020665845575
$ws.Range("J21").Value = 0.3966804807205673
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 11.1389315
$ws.Range("N21").Value = 22.277863
$ws.Range("O21").Value = 0.06331058615409579
$ws.Range("P21").Value = 0.04311697781563077
$ws.Range("Q21").Value = 1461.021356902542
$ws.Range("R21").Value = 8766.128141415255
$ws.Range("S21").Value = 0.02477356319877736
$ws.Range("T21").Value = 0.01710366348712245

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 13.634385
$ws.Range("H22").Value = 27.26877
$ws.Range("I22").Value = 0.04067566515680266
$ws.Range("J22").Value = 0.02748983251202704
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 10.99123566666667
$ws.Range("N22").Value = 32.973707
$ws.Range("O22").Value = 0.06247112414817082
$ws.Range("P22").Value = 0.06381790718517791
$ws.Range("Q22").Value = 149.858738705065
$ws.Range("R22").Value = 899.1524322303901
$ws.Range("S22").Value = 0.002541054527820045
$ws.Range("T22").Value = 0.001754343579788628

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 13.634385
$ws.Range("H23").Value = 27.26877
$ws.Range("I23").Value = 0.04067566515680266
$ws.Range("J23").Value = 0.02748983251202704
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 16.58451233333333
$ws.Range("N23").Value = 49.753537
$ws.Range("O23").Value = 0.09426175184784683
$ws.Range("P23").Value = 0.09629389277948987
$ws.Range("Q23").Value = 226.119626189915
$ws.Range("R23").Value = 1356.71775713949
$ws.Range("S23").Value = 0.003834159455256642
$ws.Range("T23").Value = 0.002647102984439266

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 13.634385
$ws.Range("H24").Value = 27.26877
$ws.Range("I24").Value = 0.04067566515680266
$ws.Range("J24").Value = 0.02748983251202704
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 78.20455633333331
$ws.Range("N24").Value = 234.613669
$ws.Range("O24").Value = 0.4444929301687811
$ws.Range("P24").Value = 0.4540755260734272
$ws.Range("Q24").Value = 1066.271029802855
$ws.Range("R24").Value = 6397.626178817129
$ws.Range("S24").Value = 0.01808004559211141
$ws.Range("T24").Value = 0.01248246015956908

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 13.634385
$ws.Range("H25").Value = 27.26877
$ws.Range("I25").Value = 0.04067566515680266
$ws.Range("J25").Value = 0.02748983251202704
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 59.02182200000001
$ws.Range("N25").Value = 177.065466
$ws.Range("O25").Value = 0.3354636076811054
$ws.Range("P25").Value = 0.3426956961462742
$ws.Range("Q25").Value = 804.7262445494702
$ws.Range("R25").Value = 4828.357467296821
$ws.Range("S25").Value = 0.01364520537832966
$ws.Range("T25").Value = 0.009420647289653588

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 13.634385
$ws.Range("H26").Value = 27.26877
$ws.Range("I26").Value = 0.04067566515680266
$ws.Range("J26").Value = 0.02748983251202704
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 11.1389315
$ws.Range("N26").Value = 22.277863
$ws.Range("O26").Value = 0.06331058615409579
$ws.Range("P26").Value = 0.04311697781563077
$ws.Range("Q26").Value = 151.8724805596275
$ws.Range("R26").Value = 607.48992223851
$ws.Range("S26").Value = 0.002575200203284907
$ws.Range("T26").Value = 0.001185278498576475

Write-Host "Updated cells for rows 2-26"
